# Auto-generated edit script applying the Halicarnassus_Profits.xlsx diff
# Updates LevePriceNQ/HQ-derived numeric columns (H-N) on specific rows
# across all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 510.6
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H15").Value = 498.6
$ws.Range("I15").Value = 498.6
$ws.Range("K15").Value = 1495.8
$ws.Range("M15").Value = -1326.8

$ws.Range("H39").Value = 179.16667
$ws.Range("I39").Value = 102.333336
$ws.Range("J39").Value = 332.83334
$ws.Range("K39").Value = 307.000008
$ws.Range("L39").Value = 998.5000200000001
$ws.Range("M39").Value = -11.00000799999998
$ws.Range("N39").Value = -1590.50002

$ws.Range("H43").Value = 985.5
$ws.Range("I43").Value = 985.5
$ws.Range("K43").Value = 985.5
$ws.Range("M43").Value = -916.5

$ws.Range("H48").Value = 1500
$ws.Range("J48").Value = 1500
$ws.Range("L48").Value = 4500
$ws.Range("N48").Value = -5084

$ws.Range("H56").Value = 1500
$ws.Range("J56").Value = 1500
$ws.Range("L56").Value = 4500
$ws.Range("N56").Value = -5568

$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594

$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596

$ws.Range("H112").Value = 1755.2222
$ws.Range("J112").Value = 1185.2858
$ws.Range("L112").Value = 3555.8574
$ws.Range("N112").Value = -5771.857400000001

$ws.Range("H113").Value = 7998.5713
$ws.Range("I113").Value = 7331.6665
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 7331.6665
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -4077.6665
$ws.Range("N113").Value = -18508

$ws.Range("H116").Value = 3612.25
$ws.Range("I116").Value = 3699.6
$ws.Range("J116").Value = 3466.6667
$ws.Range("K116").Value = 3699.6
$ws.Range("L116").Value = 3466.6667
$ws.Range("M116").Value = -257.5999999999999
$ws.Range("N116").Value = -10350.6667

$ws.Range("H137").Value = 2393.0908
$ws.Range("I137").Value = 631
$ws.Range("K137").Value = 1893
$ws.Range("M137").Value = 657

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 265.85715
$ws.Range("J4").Value = 706
$ws.Range("L4").Value = 706
$ws.Range("N4").Value = -938

$ws.Range("H5").Value = 90
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 90
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 22
$ws.Range("N5").ClearContents()

$ws.Range("H39").Value = 7779.8
$ws.Range("I39").Value = 6724.75
$ws.Range("K39").Value = 6724.75
$ws.Range("M39").Value = -6204.75

$ws.Range("H45").Value = 3862.3333
$ws.Range("I45").Value = 2218.5
$ws.Range("K45").Value = 2218.5
$ws.Range("M45").Value = -1841.5

$ws.Range("H132").Value = 4487
$ws.Range("I132").Value = 5164.4
$ws.Range("K132").Value = 15493.2
$ws.Range("M132").Value = -12963.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 90
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 25
$ws.Range("N4").ClearContents()

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H22").Value = 223.11111
$ws.Range("I22").Value = 239.875
$ws.Range("K22").Value = 239.875
$ws.Range("M22").Value = -66.875

$ws.Range("H88").Value = 16632.5
$ws.Range("J88").Value = 16632.5
$ws.Range("L88").Value = 16632.5
$ws.Range("N88").Value = -17444.5

$ws.Range("H91").Value = 16632.5
$ws.Range("J91").Value = 16632.5
$ws.Range("L91").Value = 16632.5
$ws.Range("N91").Value = -19440.5

$ws.Range("H99").Value = 1849.2
$ws.Range("I99").Value = 1311.75
$ws.Range("J99").Value = 3999
$ws.Range("K99").Value = 1311.75
$ws.Range("L99").Value = 3999
$ws.Range("M99").Value = 186.25
$ws.Range("N99").Value = -6995

$ws.Range("H105").Value = 1825.8334
$ws.Range("I105").Value = 1669.3334
$ws.Range("J105").Value = 1982.3334
$ws.Range("K105").Value = 1669.3334
$ws.Range("L105").Value = 1982.3334
$ws.Range("M105").Value = 77.66660000000002
$ws.Range("N105").Value = -5476.3334

$ws.Range("H134").Value = 1299.3334
$ws.Range("I134").Value = 1299.3334
$ws.Range("K134").Value = 3898.0002
$ws.Range("M134").Value = -1363.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H22").Value = 776.2308
$ws.Range("I22").Value = 779.4
$ws.Range("J22").Value = 765.6667
$ws.Range("K22").Value = 779.4
$ws.Range("L22").Value = 765.6667
$ws.Range("M22").Value = -429.4
$ws.Range("N22").Value = -1465.6667

$ws.Range("H60").Value = 60722.855
$ws.Range("J60").Value = 67452
$ws.Range("L60").Value = 67452
$ws.Range("N60").Value = -68474

$ws.Range("H68").Value = 66382.5
$ws.Range("J68").Value = 66382.5
$ws.Range("L68").Value = 66382.5
$ws.Range("N68").Value = -67880.5

$ws.Range("H71").Value = 66382.5
$ws.Range("J71").Value = 66382.5
$ws.Range("L71").Value = 199147.5
$ws.Range("N71").Value = -206635.5

$ws.Range("H107").Value = 1075.8889
$ws.Range("I107").Value = 467.4
$ws.Range("K107").Value = 467.4
$ws.Range("M107").Value = 1452.6

$ws.Range("H122").Value = 2886
$ws.Range("I122").Value = 2886
$ws.Range("K122").Value = 8658
$ws.Range("M122").Value = -6208

$ws.Range("H132").Value = 1447.5
$ws.Range("I132").Value = 1447.5
$ws.Range("K132").Value = 4342.5
$ws.Range("M132").Value = -1812.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 703333.7
$ws.Range("I4").Value = 1000000.5
$ws.Range("J4").Value = 110000
$ws.Range("K4").Value = 3000001.5
$ws.Range("L4").Value = 330000
$ws.Range("M4").Value = -2999889.5
$ws.Range("N4").Value = -330224

$ws.Range("H34").Value = 840.44446
$ws.Range("J34").Value = 1036.5714
$ws.Range("L34").Value = 3109.7142
$ws.Range("N34").Value = -3277.7142

$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H82").Value = 7000
$ws.Range("J82").Value = 7000
$ws.Range("L82").Value = 21000
$ws.Range("N82").Value = -21812

$ws.Range("H85").Value = 7000
$ws.Range("J85").Value = 7000
$ws.Range("L85").Value = 21000
$ws.Range("N85").Value = -23808

$ws.Range("H140").Value = 1907.8125
$ws.Range("I140").Value = 1701.6666
$ws.Range("K140").Value = 5104.9998
$ws.Range("M140").Value = 75.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1999.6666
$ws.Range("I31").Value = 1999.6666
$ws.Range("K31").Value = 1999.6666
$ws.Range("M31").Value = -1707.6666

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H37").Value = 1999.6666
$ws.Range("I37").Value = 1999.6666
$ws.Range("K37").Value = 1999.6666
$ws.Range("M37").Value = -1722.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2019.4
$ws.Range("I10").Value = 2366
$ws.Range("J10").Value = 1499.5
$ws.Range("K10").Value = 2366
$ws.Range("L10").Value = 1499.5
$ws.Range("M10").Value = -2226
$ws.Range("N10").Value = -1779.5

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 4999.5
$ws.Range("I122").Value = 4999
$ws.Range("K122").Value = 14997
$ws.Range("M122").Value = -12547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()

$ws.Range("H58").Value = 11500
$ws.Range("I58").Value = 11500
$ws.Range("K58").Value = 11500
$ws.Range("M58").Value = -11192

$ws.Range("H113").Value = 698.3
$ws.Range("I113").Value = 585.5
$ws.Range("K113").Value = 1756.5
$ws.Range("M113").Value = 413.5

$ws.Range("H122").Value = 1087
$ws.Range("I122").Value = 625.75
$ws.Range("K122").Value = 1877.25
$ws.Range("M122").Value = 572.75
